$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds values that look numeric (e.g. "6.90", "42.056.90")
# but are authored as literal text in the source data (note the double-dot
# thousands separators, and preserved trailing zeros). Force each target cell to
# Text format first so Excel does not reinterpret/round them as numbers.
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D9", "D10", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D21", "D22", "D23", "D26", "D30", "D31", "D32", "D33", "D34", "D39", "D41", "D42", "D43", "D44", "D45", "D49", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "41.979.20"
$ws.Range("E2").Value = "  -0.98%  "
$ws.Range("D3").Value = "2.237.94"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "242.44"
$ws.Range("E5").Value = "  -1.08%  "
$ws.Range("D6").Value = "0.624"
$ws.Range("E6").Value = "  -0.77%  "
$ws.Range("D7").Value = "74.03"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "0.596"
$ws.Range("E9").Value = "  -3.69%  "
$ws.Range("D10").Value = "42.05"
$ws.Range("E10").Value = "  -2.31%  "
$ws.Range("E11").Value = "  -1.67%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.103"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "6.90"
$ws.Range("E13").Value = "  -2.98%  "
$ws.Range("D14").Value = "2.571.83"
$ws.Range("E14").Value = "  +0.32%  "
$ws.Range("D15").Value = "14.32"
$ws.Range("E15").Value = "  -0.79%  "
$ws.Range("D16").Value = "0.836"
$ws.Range("E16").Value = "  -1.68%  "
$ws.Range("D17").Value = "2.232.39"
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").Value = "41.894.58"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("E19").Value = "  -5.79%  "
$ws.Range("E20").Value = "  +0.65%  "
$ws.Range("D21").Value = "72.49"
$ws.Range("D22").Value = "11.05"
$ws.Range("E22").Value = "  +6.50%  "
$ws.Range("D23").Value = "229.61"
$ws.Range("E23").Value = "  -0.61%  "
$ws.Range("E24").Value = "  -6.06%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "11.31"
$ws.Range("E26").Value = "  -3.44%  "
$ws.Range("E27").Value = "  -1.59%  "
$ws.Range("E28").Value = "  -0.91%  "
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("D30").Value = "167.47"
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("D31").Value = "20.56"
$ws.Range("E31").Value = "  -1.71%  "
$ws.Range("D32").Value = "5.57"
$ws.Range("E32").Value = "  -4.78%  "
$ws.Range("D33").Value = "0.0797"
$ws.Range("E33").Value = "  -0.93%  "
$ws.Range("D34").Value = "30.37"
$ws.Range("E34").Value = "  +2.48%  "
$ws.Range("E35").Value = "  -0.63%  "
$ws.Range("E36").Value = "  -6.87%  "
$ws.Range("E37").Value = "  -3.67%  "
$ws.Range("E38").Value = "  -1.24%  "
$ws.Range("D39").Value = "13.05"
$ws.Range("E39").Value = "  -0.94%  "
$ws.Range("E40").Value = "  -1.88%  "
$ws.Range("D41").Value = "5.66"
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("D42").Value = "64.31"
$ws.Range("E42").Value = "  +1.83%  "
$ws.Range("D43").Value = "0.197"
$ws.Range("E43").Value = "  -1.79%  "
$ws.Range("D44").Value = "8.69"
$ws.Range("E44").Value = "  -1.42%  "
$ws.Range("D45").Value = "103.07"
$ws.Range("E45").Value = "  -2.10%  "
$ws.Range("E46").Value = "  -1.94%  "
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("D49").Value = "2.32"
$ws.Range("E49").Value = "  -2.03%  "
$ws.Range("E50").Value = "  -1.49%  "
$ws.Range("D51").Value = "2.447.06"
$ws.Range("E51").Value = "  +0.04%  "
